$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 211 (open/high/low/close values revised) ---
$ws.Range("C211").Value = 9409316708800
$ws.Range("D211").Value = 9409316708800
$ws.Range("E211").Value = 9409316708800
$ws.Range("F211").Value = 9409316708800

# --- Update existing row 213 (open/high/low/close values revised) ---
$ws.Range("C213").Value = 11478475285200
$ws.Range("D213").Value = 11478475285200
$ws.Range("E213").Value = 11478475285200
$ws.Range("F213").Value = 11478475285200

# --- Append new rows 214-216 with the same formatting as row 213 ---
$ws.Range("A213:G213").Copy($ws.Range("A214:G214"))
$ws.Range("A213:G213").Copy($ws.Range("A215:G215"))
$ws.Range("A213:G213").Copy($ws.Range("A216:G216"))

# Row 214
$ws.Range("A214").Value = 45139.41666666666
$ws.Range("B214").Value = "ECONOMICS:TRM2"
$ws.Range("C214").Value = 12025285811100
$ws.Range("D214").Value = 12025285811100
$ws.Range("E214").Value = 12025285811100
$ws.Range("F214").Value = 12025285811100
$ws.Range("G214").Value = 0

# Row 215
$ws.Range("A215").Value = 45170.41666666666
$ws.Range("B215").Value = "ECONOMICS:TRM2"
$ws.Range("C215").Value = 12349311426800
$ws.Range("D215").Value = 12349311426800
$ws.Range("E215").Value = 12349311426800
$ws.Range("F215").Value = 12349311426800
$ws.Range("G215").Value = 0

# Row 216
$ws.Range("A216").Value = 45200.45833333334
$ws.Range("B216").Value = "ECONOMICS:TRM2"
$ws.Range("C216").Value = 12763732048500
$ws.Range("D216").Value = 12763732048500
$ws.Range("E216").Value = 12763732048500
$ws.Range("F216").Value = 12763732048500
$ws.Range("G216").Value = 0
